$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 2 (shifts existing rows 2..25 down to 3..26),
# mirroring a brand-new IPO record being pushed to the top of the dataset.
$ws.Rows.Item(2).Insert(-4121)

# The insert copies the bold/centered formatting down from the header row.
# Force the text-bearing columns to Text format before writing the values,
# so date-looking strings (and the "-" placeholders) are stored as literal
# shared strings instead of being auto-converted into date serial numbers.
$textColumns = @("A","B","C","D","E","H","K","M","N","P","Q","R","S","T")
foreach ($col in $textColumns) {
    $ws.Range($col + "2").NumberFormat = "@"
}

$ws.Range("A2").Value = "2023-11-16"
$ws.Range("B2").Value = "에이에스텍"
$ws.Range("C2").Value = "미래"
$ws.Range("D2").Value = "2023-11-21"
$ws.Range("E2").Value = "2023-11-28"
$ws.Range("F2").Value = 39396000
$ws.Range("G2").Value = 1407000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 21000
$ws.Range("J2").Value = 25000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 28000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 40.01421464108032
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "1356 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# Remove the Text number-format / bold-header residue now that the values are
# committed, so the new row matches the plain look of the other data rows.
$ws.Rows.Item(2).ClearFormats()
